$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 4825
$ws.Range("E3").Value = 16877
$ws.Range("E4").Value = 15186
$ws.Range("E5").Value = 6872
$ws.Range("E6").Value = 8510
$ws.Range("E7").Value = 12634
$ws.Range("E8").Value = 7311
$ws.Range("E9").Value = 1390
$ws.Range("E10").Value = 15772
$ws.Range("E11").Value = 10721
$ws.Range("E12").Value = 11069
$ws.Range("E13").Value = 1295
